# Apply price/volume updates to the cryptos worksheet (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.990.81"
$ws.Range("E2").Value = "  +1.59%  "

$ws.Range("D3").Value = "1.752.25"
$ws.Range("E3").Value = "  -0.54%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").Value = "'336.07"
$ws.Range("E5").Value = "  +0.15%  "

$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("D7").Value = "'0.3848"
$ws.Range("E7").Value = "  +0.41%  "

$ws.Range("D8").Value = "'0.3407"
$ws.Range("E8").Value = "  +0.15%  "

$ws.Range("D9").Value = "'45.80"
$ws.Range("E9").Value = "  -2.30%  "

$ws.Range("D10").Value = "'1.120"
$ws.Range("E10").Value = "  -1.49%  "

$ws.Range("E11").Value = "  -2.18%  "

$ws.Range("D12").Value = "'22.70"
$ws.Range("E12").Value = "  +1.61%  "

$ws.Range("D13").Value = "'1.001"
$ws.Range("E13").Value = "  -0.11%  "

$ws.Range("D14").Value = "'6.172"
$ws.Range("E14").Value = "  -2.63%  "

$ws.Range("D15").Value = "'7.120"
$ws.Range("E15").Value = "  +1.38%  "

$ws.Range("D16").Value = "1.751.69"
$ws.Range("E16").Value = "  -0.64%  "

$ws.Range("D18").Value = "'0.06618"
$ws.Range("E18").Value = "  -0.57%  "

$ws.Range("D19").Value = "'79.24"
$ws.Range("E19").Value = "  -3.56%  "

$ws.Range("D20").Value = "'0.9995"
$ws.Range("E20").Value = "  -0.17%  "

$ws.Range("D21").Value = "'16.76"
$ws.Range("E21").Value = "  -3.23%  "

$ws.Range("D22").Value = "'6.192"
$ws.Range("E22").Value = "  -3.11%  "

$ws.Range("D23").Value = "27.992.55"
$ws.Range("E23").Value = "  +1.55%  "

$ws.Range("D24").Value = "'11.67"
$ws.Range("E24").Value = "  -2.92%  "

$ws.Range("D25").Value = "'2.397"
$ws.Range("E25").Value = "  +0.73%  "

$ws.Range("D26").Value = "'153.58"
$ws.Range("E26").Value = "  +0.71%  "

$ws.Range("E27").Value = "  -3.61%  "

$ws.Range("D28").Value = "'2.301"
$ws.Range("E28").Value = "  -4.62%  "

$ws.Range("D29").Value = "1.950.77"
$ws.Range("E29").Value = "  -0.63%  "

$ws.Range("D30").Value = "'1.260"
$ws.Range("E30").Value = "  -11.76%  "

$ws.Range("D31").Value = "'131.25"
$ws.Range("E31").Value = "  -2.30%  "

$ws.Range("D32").Value = "'4.027"
$ws.Range("E32").Value = "  +1.73%  "

$ws.Range("D33").Value = "'5.847"
$ws.Range("E33").Value = "  -4.34%  "

$ws.Range("E34").Value = "  +0.34%  "

$ws.Range("D35").Value = "'12.22"
$ws.Range("E35").Value = "  -3.67%  "

$ws.Range("D36").Value = "'1.541"
$ws.Range("E36").Value = "  +2.05%  "

$ws.Range("D37").Value = "'0.6559"
$ws.Range("E37").Value = "  -3.10%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02284"
$ws.Range("E38").Value = "  -5.19%  "

$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'5.144"
$ws.Range("E39").Value = "  -3.43%  "

$ws.Range("D40").Value = "'0.06133"
$ws.Range("E40").Value = "  -2.45%  "

$ws.Range("E41").Value = "  -3.70%  "

$ws.Range("D42").Value = "'1.208"
$ws.Range("E42").Value = "  -2.98%  "

$ws.Range("D43").Value = "'8.010"
$ws.Range("E43").Value = "  -2.97%  "

$ws.Range("D44").Value = "'0.9992"
$ws.Range("E44").Value = "  -0.13%  "

$ws.Range("D45").Value = "'13.74"
$ws.Range("E45").Value = "  -3.23%  "

$ws.Range("D46").Value = "'3.838"
$ws.Range("E46").Value = "  +0.38%  "

$ws.Range("D47").Value = "'0.6044"
$ws.Range("E47").Value = "  -3.21%  "

$ws.Range("D48").Value = "'126.68"
$ws.Range("E48").Value = "  -3.52%  "

$ws.Range("D49").Value = "'2.005"
$ws.Range("E49").Value = "  -3.28%  "

$ws.Range("D50").Value = "'1.173"
$ws.Range("E50").Value = "  +2.46%  "

$ws.Range("D51").Value = "'1.108"
$ws.Range("E51").Value = "  +4.39%  "
